$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; existing data (old rows 1-384) shifts down to rows 3-386
$ws.Rows("1:2").Insert()

# New row 2: "Node1" / "Node2" column headers
$ws.Range("A2").Value = "Node1"
$ws.Range("B2").Value = "Node2"

# New row 1: "Edge List" title, merged across A1:B1
$ws.Range("A1").Value = "Edge List"
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").Merge()

# Center the Node1/Node2 header row
$ws.Range("A2:B2").HorizontalAlignment = -4108

# Update the active selection to match the post-edit state (was E6, now E8 after +2 row shift)
$ws.Range("E8").Select()

# Explicit portrait page orientation
$ws.PageSetup.Orientation = 1

Write-Output "done"
